# Add data for 2021-10-26 (refresh carjacking arrests-by-month YoY data
# through 2021-10-18, previously through 2021-10-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet / tab name reflects the new "as-of" date.
$wb.Worksheets.Item(1).Name = "Through 2021-10-18"

# --- Row 9 (July) ---
$ws.Range("U9").Value = 138
$ws.Range("V9").Value = 0.08

# --- Row 12 (October) ---
$ws.Range("A12").Value = "October (through 10-18)"
$ws.Range("C12").Value = 16
$ws.Range("D12").Value = 0.0588
$ws.Range("H12").Value = 7
$ws.Range("J12").Value = 0.2258
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.0652
$ws.Range("N12").Value = 3
$ws.Range("O12").Value = 25
$ws.Range("P12").Value = 0.1071
$ws.Range("R12").Value = 86
$ws.Range("U12").Value = 112

# --- Row 13 (Total) ---
$ws.Range("C13").Value = 212
$ws.Range("D13").Value = 0.1276
$ws.Range("H13").Value = 57
$ws.Range("J13").Value = 0.0866
$ws.Range("K13").Value = 64
$ws.Range("M13").Value = 0.1077
$ws.Range("N13").Value = 46
$ws.Range("O13").Value = 404
$ws.Range("P13").Value = 0.1022
$ws.Range("R13").Value = 934
$ws.Range("S13").Value = 0.0537
$ws.Range("U13").Value = 1277
$ws.Range("V13").Value = 0.061
